# Balance tweaks and pdf creation
$wb = $excel.ActiveWorkbook

# --- New "Economy" sheet ---
$economy = $wb.Worksheets.Add()
$economy.Name = "Economy"

# Excel inserts new sheets before the active sheet by default; move it to
# the end so the tab order is Deck, VPs, Economy. Re-fetch the worksheet
# object by name afterwards since the move changes which sheet sits at
# the position the old reference was bound to.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$economy.Move($null, $lastSheet)
$economy = $wb.Worksheets.Item("Economy")

$economy.Range("A1").Value = "Resource"
$economy.Range("C1").Value = "Total Cost"
$economy.Range("B1").Value = "# Cards Costing"

$economy.Range("A2").Value = "Wood"
$economy.Range("A3").Value = "Steel"
$economy.Range("A4").Value = "Stone"
$economy.Range("A5").Value = "Gold"

$economy.Range("B2").Formula = "=COUNTIF(Deck!D:D,Economy!A2) + COUNTIF(Deck!F:F,Economy!A2)"
$economy.Range("C2").Formula = "=SUMIF(Deck!D:D,Economy!A2,Deck!C:C) + SUMIF(Deck!F:F,Economy!A2,Deck!E:E)"

$economy.Range("B3").Formula = "=COUNTIF(Deck!D:D,Economy!A3) + COUNTIF(Deck!F:F,Economy!A3)"
$economy.Range("C3").Formula = "=SUMIF(Deck!D:D,Economy!A3,Deck!C:C) + SUMIF(Deck!F:F,Economy!A3,Deck!E:E)"

$economy.Range("B4").Formula = "=COUNTIF(Deck!D:D,Economy!A4) + COUNTIF(Deck!F:F,Economy!A4)"
$economy.Range("C4").Formula = "=SUMIF(Deck!D:D,Economy!A4,Deck!C:C) + SUMIF(Deck!F:F,Economy!A4,Deck!E:E)"

$economy.Range("B5").Formula = "=COUNTIF(Deck!D:D,Economy!A5) + COUNTIF(Deck!F:F,Economy!A5)"
$economy.Range("C5").Formula = "=SUMIF(Deck!D:D,Economy!A5,Deck!C:C) + SUMIF(Deck!F:F,Economy!A5,Deck!E:E)"

$economy.Range("A1").Font.Bold = $true
$economy.Range("B1").Font.Bold = $true
$economy.Range("C1").Font.Bold = $true
$economy.Range("B1").HorizontalAlignment = -4108
$economy.Range("B2:B5").HorizontalAlignment = -4108
$economy.Range("C2:C5").HorizontalAlignment = -4108
$economy.Range("B2:C5").Font.Bold = $false

$economy.Range("B3").Select()

# --- Deck sheet balance tweak: row 13 "Lumberjack" -> "Ladder" ---
$deck = $wb.Worksheets.Item("Deck")
$deck.Range("A13").Value = "Ladder"
$deck.Range("D13").Value = "Wood"
$deck.Range("H13").Value = "hole-ladder"
$deck.Range("I13").Value = "Draw a card from your deck, if it's Wood, draw two more cards."

$wb.Worksheets.Item("Deck").Activate()
# Move the visible selection to match the edited cell
$deck.Range("H13").Select()
